$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "11/04/2022"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "12:51"
$ws.Range("C6").Value = "12:53"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = "EURUSD"
$ws.Range("J6").Value = -12.22
$ws.Range("K6").Value = "LOSS"

# Row 7
$ws.Range("A7").Value = "13/04/2022"
$ws.Range("B7").Value = "14:15"
$ws.Range("C7").Value = "14:17"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = "EURUSD"
$ws.Range("J7").Value = -14
$ws.Range("K7").Value = "LOSS"
